{"js": "// Replace the date line and each \"NNN\u00d7N=\" drill cell with its new value.\n// Each old string is unique in the document, so a plain body.search() +\n// insertText(\"Replace\") round-trip is safe for every entry.\nconst replacements = [\n  [\"2026-01-09 Friday\", \"2026-01-10 Saturday\"],\n  [\"660\u00d78=\", \"479\u00d77=\"],\n  [\"375\u00d73=\", \"187\u00d77=\"],\n  [\"931\u00d76=\", \"294\u00d75=\"],\n  [\"913\u00d75=\", \"754\u00d79=\"],\n  [\"565\u00d79=\", \"928\u00d72=\"],\n  [\"839\u00d72=\", \"897\u00d77=\"],\n  [\"368\u00d78=\", \"516\u00d78=\"],\n  [\"588\u00d75=\", \"977\u00d75=\"],\n  [\"816\u00d78=\", \"137\u00d76=\"],\n  [\"507\u00d74=\", \"677\u00d75=\"],\n  [\"869\u00d76=\", \"152\u00d74=\"],\n  [\"807\u00d73=\", \"525\u00d76=\"],\n  [\"215\u00d77=\", \"305\u00d74=\"],\n  [\"832\u00d74=\", \"687\u00d78=\"],\n  [\"316\u00d79=\", \"683\u00d74=\"],\n  [\"817\u00d76=\", \"864\u00d76=\"],\n  [\"772\u00d75=\", \"714\u00d74=\"],\n  [\"549\u00d79=\", \"907\u00d76=\"],\n  [\"748\u00d76=\", \"891\u00d78=\"],\n  [\"288\u00d74=\", \"469\u00d72=\"],\n  [\"341\u00d78=\", \"875\u00d76=\"],\n  [\"977\u00d79=\", \"852\u00d74=\"],\n  [\"254\u00d75=\", \"916\u00d77=\"],\n  [\"548\u00d75=\", \"885\u00d74=\"],\n  [\"679\u00d76=\", \"731\u00d75=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and each \"NNN\u00d7N=\" drill cell with its new value.\n# Each old string occurs exactly once in the document, so a simple\n# Find/Replace (wdReplaceAll = 2, but only ever one hit per string) scoped\n# to the whole document body is safe for every entry.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2026-01-09 Friday\", \"2026-01-10 Saturday\"),\n    @(\"660\u00d78=\", \"479\u00d77=\"),\n    @(\"375\u00d73=\", \"187\u00d77=\"),\n    @(\"931\u00d76=\", \"294\u00d75=\"),\n    @(\"913\u00d75=\", \"754\u00d79=\"),\n    @(\"565\u00d79=\", \"928\u00d72=\"),\n    @(\"839\u00d72=\", \"897\u00d77=\"),\n    @(\"368\u00d78=\", \"516\u00d78=\"),\n    @(\"588\u00d75=\", \"977\u00d75=\"),\n    @(\"816\u00d78=\", \"137\u00d76=\"),\n    @(\"507\u00d74=\", \"677\u00d75=\"),\n    @(\"869\u00d76=\", \"152\u00d74=\"),\n    @(\"807\u00d73=\", \"525\u00d76=\"),\n    @(\"215\u00d77=\", \"305\u00d74=\"),\n    @(\"832\u00d74=\", \"687\u00d78=\"),\n    @(\"316\u00d79=\", \"683\u00d74=\"),\n    @(\"817\u00d76=\", \"864\u00d76=\"),\n    @(\"772\u00d75=\", \"714\u00d74=\"),\n    @(\"549\u00d79=\", \"907\u00d76=\"),\n    @(\"748\u00d76=\", \"891\u00d78=\"),\n    @(\"288\u00d74=\", \"469\u00d72=\"),\n    @(\"341\u00d78=\", \"875\u00d76=\"),\n    @(\"977\u00d79=\", \"852\u00d74=\"),\n    @(\"254\u00d75=\", \"916\u00d77=\"),\n    @(\"548\u00d75=\", \"885\u00d74=\"),\n    @(\"679\u00d76=\", \"731\u00d75=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
